$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append a trailing space run to the end of the "Features that will be
#    tested..." paragraph.
# ---------------------------------------------------------------------------
$endRng = $d.Content
$endRng.Find.Execute("harmony together.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endRng.Collapse(0)
$endRng.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 2) Replace the old manually-numbered "8.0 Features Not to Be Tested"
#    heading paragraph with:
#      - a Heading1 paragraph that uses the automatic numbering (numId 7)
#        and carries bookmark _Toc40375384 (id 40), text "Features Not to
#        Be Tested"
#      - a brand new plain paragraph with the "not tested" body text
# ---------------------------------------------------------------------------
$headingRng = $d.Content
$headingRng.Find.Execute("8.0 Features Not to Be Tested", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $headingRng.Paragraphs(1).Range

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Heading1"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr>
<w:jc w:val="both"/>
<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>
</w:pPr>
<w:bookmarkStart w:id="40" w:name="_Toc40375384"/>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Features Not to Be Tested</w:t></w:r>
<w:bookmarkEnd w:id="40"/>
</w:p>
<w:p>
<w:r><w:t>Features that will not be tested will be the volume settings and music settings as in the hierarchy on priority they are low on the hierarchy.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$headingPara.InsertXML($xmlSnippet)

# ---------------------------------------------------------------------------
# 3) Mark a lastRenderedPageBreak right before the "9" in the
#    "9.0 Resources/Roles & Responsibilities" heading (now pushed onto the
#    new page by the extra content above).
# ---------------------------------------------------------------------------
$nineRng = $d.Content
$nineRng.Find.Execute("9.0 Resources", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nineStart = $d.Range($nineRng.Start, $nineRng.Start)

$pageBreakXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$nineStart.InsertXML($pageBreakXml)
